# Fruta / hortaliza, semanal
# New weekly price record is added at the top of the data table (row 9,
# right after the 7 "promedio nacional" style rows already present),
# pushing every existing record (previously rows 9-105) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9; everything below (old rows 9..105)
# shifts down to rows 10..106, growing the used range to A1:R106.
$ws.Rows.Item(9).EntireRow.Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 45168
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 100112022
$ws.Range("G9").Value = "Arveja Verde"
$ws.Range("H9").Value = "Perfection"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 21000
$ws.Range("L9").Value = 23000
$ws.Range("M9").Value = 22000
$ws.Range("N9").Value = "`$/malla 25 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 880
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
